# Added 4 video tutorials on Git and Github
# Sheet3 ("Git Github" tab) gets an author byline plus four new rows of
# video-tutorial links (topic name + hyperlinked Google-Drive link).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- New tutorial rows (topic name in column A, hyperlink (URL shown as
#     the cell text) in column C) ---
$ws.Range("A5").Value = "Introduction to Git & Github"
$ws.Range("C5").Value = "https://drive.google.com/file/d/1OzNETc6oe5tDewQivYgt0elq6NfQ0BRz/view?usp=sharing"
$ws.Hyperlinks.Add($ws.Range("C5"), "https://drive.google.com/file/d/1OzNETc6oe5tDewQivYgt0elq6NfQ0BRz/view?usp=sharing")

$ws.Range("A6").Value = "git status command"
$ws.Range("C6").Value = "https://drive.google.com/file/d/1RRFCBk8WkjKQj91gcCRLCrFmAXOpLVo0/view?usp=sharing"
$ws.Hyperlinks.Add($ws.Range("C6"), "https://drive.google.com/file/d/1RRFCBk8WkjKQj91gcCRLCrFmAXOpLVo0/view?usp=sharing")

$ws.Range("A7").Value = "git add command"
$ws.Range("C7").Value = "https://drive.google.com/file/d/1yBP7n8PpCNyAX4891f0uQ8FIk9brSEDg/view?usp=sharing"
$ws.Hyperlinks.Add($ws.Range("C7"), "https://drive.google.com/file/d/1yBP7n8PpCNyAX4891f0uQ8FIk9brSEDg/view?usp=sharing")

$ws.Range("A8").Value = "git init command"
$ws.Range("C8").Value = "https://drive.google.com/file/d/1NuG7XGedxh919clWwoaQYo33TTC0WhZ_/view"
$ws.Hyperlinks.Add($ws.Range("C8"), "https://drive.google.com/file/d/1NuG7XGedxh919clWwoaQYo33TTC0WhZ_/view")

# --- Byline above the table: author name + sheet heading ---
$ws.Range("A2").Value = "Utkarsh Rai"
$ws.Range("A2").Font.Color = 0

$ws.Range("A1").Value = "Git Github"
$ws.Range("A1").Font.Color = 0
$ws.Range("A1").Interior.Color = 16777215

# --- Column widths for the new content (27.43 / 82.43 characters; values
#     chosen so the engine's pixel-grid rounding lands on the closest
#     achievable width) ---
$ws.Columns.Item(1).ColumnWidth = 26.65
$ws.Columns.Item(3).ColumnWidth = 81.65

# --- Make this sheet the active / selected one, with C8 selected ---
$ws.Activate()
$ws.Range("C8").Select()
